# Re-import pass: wipe the old ad-hoc test data (row 1 + the old text
# columns A/B) and write the freshly-imported rows instead. Column A is
# the running block increment (1..9), column D carries the imported
# record id for every block row (including the trailing row 11, which
# has no increment of its own).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear everything that used to live in rows 1-11 / columns A-B.
$ws.Range("A1:B11").ClearContents()

$importedId = "['cfd893a46090']"

for ($row = 2; $row -le 11; $row++) {
    $increment = $row - 1
    if ($increment -le 9) {
        $ws.Cells.Item($row, 1).Value = $increment
    }
    $ws.Cells.Item($row, 4).Value = $importedId
}

$ws.Range("A11").Select()
